# Auto-generated edit script: updates crypto price/volume table cells
# to match the "Updated cryptos list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.959.61"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").Value = "1.818.32"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Formula = "'309.95"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").Formula = "'0.4666"
$ws.Range("E7").Value = "  +0.63%  "
$ws.Range("E8").Value = "  -0.88%  "
$ws.Range("D9").Formula = "'0.07357"
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").Formula = "'0.8724"
$ws.Range("E10").Value = "  -0.33%  "
$ws.Range("D11").Formula = "'20.26"
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("D12").Value = "1.835.80"
$ws.Range("E12").Value = "  +1.85%  "
$ws.Range("D13").Formula = "'5.411"
$ws.Range("E13").Value = "  +0.95%  "
$ws.Range("D14").Formula = "'0.07113"
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("D15").Formula = "'6.512"
$ws.Range("E15").Value = "  +0.11%  "
$ws.Range("D16").Formula = "'91.54"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").Formula = "'14.64"
$ws.Range("E20").Value = "  -0.63%  "
$ws.Range("D21").Value = "26.978.53"
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("D22").Formula = "'5.295"
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("D23").Formula = "'10.60"
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("D24").Value = "2.050.80"
$ws.Range("E24").Value = "  +1.68%  "
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("D26").Formula = "'150.85"
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("E28").Value = "  +0.35%  "
$ws.Range("E29").Value = "  -0.95%  "
$ws.Range("D30").Formula = "'117.21"
$ws.Range("E30").Value = "  +1.11%  "
$ws.Range("D31").Formula = "'0.08899"
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("D32").Formula = "'0.7589"
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("E33").Value = "  +0.90%  "
$ws.Range("E34").Value = "  +0.87%  "
$ws.Range("D35").Formula = "'2.912"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("E37").Value = "  -0.54%  "
$ws.Range("D38").Formula = "'0.05301"
$ws.Range("E38").Value = "  +0.84%  "
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("D40").Formula = "'2.971"
$ws.Range("E40").Value = "  +1.79%  "
$ws.Range("D41").Formula = "'7.200"
$ws.Range("E41").Value = "  +0.47%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Formula = "'2.373"
$ws.Range("E42").Value = "  -2.44%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Formula = "'0.5299"
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").Formula = "'0.1651"
$ws.Range("E44").Value = "  -0.45%  "
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("E46").Value = "  -1.29%  "
$ws.Range("D47").Formula = "'10.50"
$ws.Range("E47").Value = "  +1.43%  "
$ws.Range("E48").Value = "  +0.15%  "
$ws.Range("D49").Formula = "'103.55"
$ws.Range("E49").Value = "  +0.65%  "
$ws.Range("D50").Formula = "'1.666"
$ws.Range("E50").Value = "  -0.39%  "
$ws.Range("D51").Formula = "'0.06296"
$ws.Range("E51").Value = "  +0.13%  "
